$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column (D2:D51) as Text before writing, so that
# numeric-looking strings (e.g. "6.48") are preserved verbatim as text
# instead of being coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '59.598.46'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '2.645.41'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '518.68'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').Value = '147.51'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').Value = '0.576'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = '2.672.27'
$ws.Range('E9').Value = '  +2.06%  '
$ws.Range('D10').Value = '6.48'
$ws.Range('E10').Value = '  +3.03%  '
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('D14').Value = '3.112.83'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').Value = '59.520.01'
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('D16').Value = '21.33'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').Value = '2.665.16'
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('D19').Value = '4.63'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '346.97'
$ws.Range('E20').Value = '  +1.27%  '
$ws.Range('D21').Value = '10.59'
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('D22').Value = '6.23'
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').Value = '61.40'
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('D26').Value = '2.764.23'
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.993'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = '0.162'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '0.0₃0825'
$ws.Range('E29').Value = '  +1.95%  '
$ws.Range('D30').Value = '7.20'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('D32').Value = '6.55'
$ws.Range('E32').Value = '  +9.11%  '
$ws.Range('D33').Value = '19.15'
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').Value = '1.59'
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').Value = '149.89'
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('E36').Value = '  +15.06%  '
$ws.Range('E37').Value = '  +2.83%  '
$ws.Range('D38').Value = '1.18'
$ws.Range('E38').Value = '  +3.39%  '
$ws.Range('D39').Value = '0.877'
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('D40').Value = '36.59'
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('D41').Value = '3.75'
$ws.Range('E41').Value = '  +3.46%  '
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '285.97'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('D45').Value = '0.0998'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('D46').Value = '0.994'
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('D47').Value = '19.85'
$ws.Range('E47').Value = '  +2.16%  '
$ws.Range('D48').Value = '0.0548'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('E50').Value = '  +1.37%  '
$ws.Range('E51').Value = '  -1.27%  '

# Restore the default (style-less) cell style on the Price column now
# that the values are safely stored as text, so no stray number format
# is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
